# Automatic update of files.
# Bump the "Förändrad" (last changed) date in column C for all data rows
# (rows 2 through 103) from 2023-09-03 (serial 45172) to 2023-09-06 (serial 45175).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 103; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value = 45175
    }
}
